$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New animal-image rows appended to the table (IA_0023 .. IA_0029)
$rows = @(
    @{ Row=24; A="IA_0023"; B="IMG_0332.jpeg"; D="Daniela Subotic"; E="CC BY 4.0"; F="data/Multimedia_Data/Image_Animal/"; G="Kenzi and Loki sleeping in the office";           H="A_004, A_005"; I=20 },
    @{ Row=25; A="IA_0024"; B="IMG_0837.jpeg"; D="Daniela Subotic"; E="CC BY 4.0"; F="data/Multimedia_Data/Image_Animal/"; G="Kenzi and Loki sleeping in the office part II";   H="A_004, A_005"; I=21 },
    @{ Row=26; A="IA_0025"; B="IMG_1226.jpeg"; D="Daniela Subotic"; E="CC BY 4.0"; F="data/Multimedia_Data/Image_Animal/"; G="Kenzi and Loki sleeping in the office part III";  H="A_004, A_005"; I=22 },
    @{ Row=27; A="IA_0026"; B="IMG_1659.jpeg"; D="Daniela Subotic"; E="CC BY 4.0"; F="data/Multimedia_Data/Image_Animal/"; G="Atli sleeping during a meeting";                  H="A_001";        I=23 },
    @{ Row=28; A="IA_0027"; B="IMG_1853.jpeg"; D="Daniela Subotic"; E="CC BY 4.0"; F="data/Multimedia_Data/Image_Animal/"; G="Atli sleeping during a meeting part II";          H="A_001";        I=24 },
    @{ Row=29; A="IA_0028"; B="IMG_2164.jpeg"; D="Daniela Subotic"; E="CC BY 4.0"; F="data/Multimedia_Data/Image_Animal/"; G="Loki sleeping in the office";                     H="A_005";        I=25 },
    @{ Row=30; A="IA_0029"; B="IMG_2163.jpeg"; D="Daniela Subotic"; E="CC BY 4.0"; F="data/Multimedia_Data/Image_Animal/"; G="Kenzi sleeping in the office";                    H="A_004";        I=26 }
)

# Column A (IDs), top to bottom
foreach ($r in $rows) {
    $ws.Range("A" + $r.Row).Value = $r.A
}

# Columns B (file name) and G (label), row by row
foreach ($r in $rows) {
    $ws.Range("B" + $r.Row).Value = $r.B
    $ws.Range("G" + $r.Row).Value = $r.G
}

# Column H (animal-character id): solo Kenzi row first, then the shared
# Kenzi+Loki rows, then the solo Loki row, matching how the labels were
# cross-referenced against the character legend
$ws.Range("H30").Value = "A_004"
$ws.Range("H24").Value = "A_004, A_005"
$ws.Range("H29").Value = "A_005"
$ws.Range("H25").Value = "A_004, A_005"
$ws.Range("H26").Value = "A_004, A_005"
$ws.Range("H27").Value = "A_001"
$ws.Range("H28").Value = "A_001"

# Remaining columns: Copyright, License List, Image Directory, Seqnum
foreach ($r in $rows) {
    $ws.Range("D" + $r.Row).Value = $r.D
    $ws.Range("E" + $r.Row).Value = $r.E
    $ws.Range("F" + $r.Row).Value = $r.F
    $ws.Range("I" + $r.Row).Value = $r.I
}

# Column G got wider to fit the longer labels now in it
$ws.Columns.Item(7).ColumnWidth = 39.67

# Leave the selection where the author ended up after entering the data
$ws.Range("F28").Select() | Out-Null
